$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01037833333333333
$ws.Range("H2").Value = 0.031135
$ws.Range("I2").Value = 0.02114284782989566
$ws.Range("J2").Value = 0.02114284782989566
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 0.1486372549783334
$ws.Range("R2").Value = 1.337735294805
$ws.Range("S2").Value = 0.006236229226847496
$ws.Range("T2").Value = 0.006236229226847496
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01037833333333333
$ws.Range("H3").Value = 0.031135
$ws.Range("I3").Value = 0.02114284782989566
$ws.Range("J3").Value = 0.02114284782989566
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 0.2810885339383334
$ws.Range("R3").Value = 2.529796805445001
$ws.Range("S3").Value = 0.0117933591476341
$ws.Range("T3").Value = 0.0117933591476341
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01037833333333333
$ws.Range("H4").Value = 0.031135
$ws.Range("I4").Value = 0.02114284782989566
$ws.Range("J4").Value = 0.02114284782989566
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 0.07420290734277778
$ws.Range("R4").Value = 0.667826166085
$ws.Range("S4").Value = 0.003113259455414065
$ws.Range("T4").Value = 0.003113259455414064
# Row 5
$ws.Range("G5").Value = 0.4265683333333333
$ws.Range("I5").Value = 0.8690094132698448
$ws.Range("J5").Value = 0.8690094132698448
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 6.109260908368332
$ws.Range("R5").Value = 54.98334817531499
$ws.Range("S5").Value = 0.2563203379715071
$ws.Range("T5").Value = 0.2563203379715071
# Row 6
$ws.Range("G6").Value = 0.4265683333333333
$ws.Range("I6").Value = 0.8690094132698448
$ws.Range("J6").Value = 0.8690094132698448
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("S6").Value = 0.4847284621173306
$ws.Range("T6").Value = 0.4847284621173306
# Row 7
$ws.Range("G7").Value = 0.4265683333333333
$ws.Range("I7").Value = 0.8690094132698448
$ws.Range("J7").Value = 0.8690094132698448
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 3.049874146172777
$ws.Range("R7").Value = 27.448867315555
$ws.Range("S7").Value = 0.1279606131810071
$ws.Range("T7").Value = 0.127960613181007
# Row 8
$ws.Range("G8").Value = 0.05392066666666667
$ws.Range("H8").Value = 0.161762
$ws.Range("I8").Value = 0.1098477389002595
$ws.Range("J8").Value = 0.1098477389002595
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 0.7722453714406666
$ws.Range("R8").Value = 6.950208342966
$ws.Range("S8").Value = 0.03240035047995197
$ws.Range("T8").Value = 0.03240035047995197
# Row 9
$ws.Range("G9").Value = 0.05392066666666667
$ws.Range("H9").Value = 0.161762
$ws.Range("I9").Value = 0.1098477389002595
$ws.Range("J9").Value = 0.1098477389002595
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 1.460396448592667
$ws.Range("R9").Value = 13.143568037334
$ws.Range("S9").Value = 0.06127243817053433
$ws.Range("T9").Value = 0.06127243817053433
# Row 10
$ws.Range("G10").Value = 0.05392066666666667
$ws.Range("H10").Value = 0.161762
$ws.Range("I10").Value = 0.1098477389002595
$ws.Range("J10").Value = 0.1098477389002595
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 0.3855214613002222
$ws.Range("R10").Value = 3.469693151702
$ws.Range("S10").Value = 0.01617495024977324
$ws.Range("T10").Value = 0.01617495024977324
